# "add more input checks" - two more quiz submissions were tested, so two
# more rows of results were recorded under the header row (with an extra
# ":" added to the "Age" label), pushing all previously recorded rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Make room for the two new result rows right under the header.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# The "Age" header gained a trailing colon to match the other labels.
$ws.Cells.Item(1, 3).Value = "Возраст:"

# New submission (row 3 first, then row 2 - matches entry order).
$ws.Cells.Item(3, 1).Value = "dasd"
$ws.Cells.Item(3, 2).Value = "dasd"
$ws.Cells.Item(3, 3).Value = "dasd"
$ws.Cells.Item(3, 4).Value = "0/5"

$ws.Cells.Item(2, 1).Value = "das"
$ws.Cells.Item(2, 2).Value = "ads"
$ws.Cells.Item(2, 3).Value = "14"
$ws.Cells.Item(2, 4).Value = "0/5"
